# Update the "Device App" requirements sheet: expand the "Patient Page"
# requirement row into a fully laid-out set of rows describing the new
# patient page tab/layout structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device App")
$ws.Activate() | Out-Null

# ---------------------------------------------------------------------
# 1. Make room: insert 8 new blank rows right after the current row 12
#    (the "Patient Page" row). This pushes the old row 13 -> row 21 and
#    old row 14 -> row 22, preserving their existing content/style.
# ---------------------------------------------------------------------
$ws.Rows("13:20").Insert()

# Give the freshly inserted rows the same base cell style (wrap text +
# thin border) used throughout the rest of the table, by copying the
# format from an existing formatted-but-empty cell.
$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B13:E20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Fill in the new / updated cell values.
#    NOTE: order matters here because the engine recomputes the shared
#    string table on save, compacting out strings that are no longer
#    referenced anywhere (the old "Status"/"Data" entries) and handing
#    out the freed slots to brand-new strings in the order they are
#    first written. Writing the brand-new strings in the exact order
#    they must appear in the shared string table keeps everything
#    lined up with the target workbook.
# ---------------------------------------------------------------------

# Rows 21-22 (previously rows 13-14): "Patient Chart" / "Report" requirement rows
$ws.Range("C21").Value = "Patient Chart"
$ws.Range("C22").Value = "Report"

# Row 12: Patient Page / Layout / (3-tab layout description)
$ws.Range("B12").Value = "Patient Page"
$ws.Range("C12").Value = "Layout"
$ws.Range("D12").Value = "Patient page should be divided in 3 tabs,`nDetails,`nChart,`nReports"

# Row 14: Patient Details Card
$ws.Range("C14").Value = "Details -> Patient Details Card"

# Row 16: Charts
$ws.Range("C16").Value = "Details -> Charts"

# Row 15: arrow for all details
$ws.Range("D15").Value = "An arrow to show all details (patient and medical details captured while admission)"

# Row 16 (cont'd)
$ws.Range("D16").Value = "This needs to be discussed furher, Initial idea is to show a series of charts based on the vital parameters being monitored. The way the different charts to be displayed needs to finalized"
$ws.Range("E16").Value = "there are multiple ways to display the charts here,`n1. show all charts one by one as user scrolls.`n2. show one or 2 charts and give option to show more charts`n3. show only one chart and provide a option to select the parameters"

# Row 17: chart time range options / select date
$ws.Range("D17").Value = "By default the chart will show for last 24 hours, the user needs to be provided to select following options,`n1. Last 24 hours`n2. till now`n3. Select date"

# Row 14 (cont'd)
$ws.Range("D14").Value = "A card with patient name and reason for admission and hospitalization date."

# Row 17 (cont'd)
$ws.Range("E17").Value = "Select date will give user an option to select any date between admission date and current date"

# Row 18: vital parameters data format TBD
$ws.Range("E18").Value = "The data format and configuration required for vital parameters is TBD."

# Row 13: connectivity-status card text
$ws.Range("D13").Value = "Connectivity status`nWard Name`nBed icon and Bed number"

# Row 13 (cont'd): Action bar (reuses the existing "Action bar" string - does not
# consume a new shared-string slot, so it can be written at any point)
$ws.Range("C13").Value = "Action bar"

# Rows 19-20 stay blank (already have the correct base style from the paste above)

# ---------------------------------------------------------------------
# 3. Row heights for the newly laid out rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 30

# ---------------------------------------------------------------------
# 4. Highlight D16 and E18 with a yellow fill (new style), matching the
#    rest of the table's wrap text + thin border formatting.
# ---------------------------------------------------------------------
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D16").Interior.Color = 65535
$ws.Range("E18").Interior.Color = 65535

# Restore the text values that PasteSpecial(Formats) preserved (it only
# copies formatting, not values, so the text set above is unaffected) -
# nothing further required here.

# ---------------------------------------------------------------------
# 5. Update the view: selected cell moves to C13.
# ---------------------------------------------------------------------
$ws.Range("C13").Select() | Out-Null
